$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The diff inserts one new weekly price record for "Betarraga" at
# Vega Monumental Concepción, pushing the existing rows 213-224 down to
# 214-225 (dimension grows from A1:R224 to A1:R225). Inserting a whole
# row at 213 reproduces that shift (and copies row formatting, matching
# the preserved s="2" style on column D).
$ws.Rows(213).Insert()

# Populate the newly inserted row 213 with the new record's data.
$ws.Range("A213").Value = 11
$ws.Range("B213").Value = "Vega Monumental Concepción"
$ws.Range("C213").Value = "Bíobío"
$ws.Range("D213").Value = 44516
$ws.Range("E213").Value = 8
$ws.Range("F213").Value = 100114014
$ws.Range("G213").Value = "Betarraga"
$ws.Range("H213").Value = "Sin especificar"
$ws.Range("I213").Value = "Primera"
$ws.Range("J213").Value = 1200
$ws.Range("K213").Value = 600
$ws.Range("L213").Value = 650
$ws.Range("M213").Value = 625
$ws.Range("N213").Value = '$/paquete 5 unidades'
$ws.Range("O213").Value = "Región Metropolitana"
$ws.Range("P213").Value = 125
$ws.Range("Q213").Value = 5
$ws.Range("R213").Value = "Hortaliza"
